# Re-upload of data/type6.xlsx:
#  - column A values on rows 2-15 change from 9 to 6
#  - the sheet selection narrows from A2:C15 to A2:A15 (anchor stays A2)
#  - the workbook window was resized (32000x13660 -> 30240x11500 twips)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A data values (rows 2-15): 9 -> 6 ---
$dataRange = $ws.Range("A2:A15")
for ($i = 1; $i -le $dataRange.Rows.Count; $i++) {
    $dataRange.Cells.Item($i, 1).Value = 6
}

# --- Narrow the active selection to A2:A15 (was A2:C15) ---
$ws.Range("A2:A15").Select()

# --- Resize the workbook window to match the re-saved file ---
$excel.ActiveWindow.Width = 30240
$excel.ActiveWindow.Height = 11500
